$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C ("Dc (cm^2/s)") for "SOC" and "Initial SOC"
$ws.Columns("C:D").Insert()

# Insert one new column before the (now shifted) column F ("P") for "Dt* (cm^2/s)",
# which sits right after "Dc (cm^2/s)" (now column E)
$ws.Columns("F:F").Insert()

# Set the new header labels
$ws.Range("C1").Value = "SOC"
$ws.Range("D1").Value = "Initial SOC"
$ws.Range("F1").Value = "Dt* (cm^2/s)"

# Fill in the new SOC / Initial SOC / Dt* (cm^2/s) data for each data row
$ws.Cells.Item(2, 3).Value = 0.9709798125603629
$ws.Cells.Item(2, 4).Value = 0.9709816968840589
$ws.Cells.Item(3, 3).Value = 0.9709586836957397
$ws.Cells.Item(3, 4).Value = 0.9709779533835073
$ws.Cells.Item(4, 3).Value = 0.9626268557354161
$ws.Cells.Item(4, 4).Value = 0.9709394399062984
$ws.Cells.Item(4, 6).Value = [double]"6.367216926490938e-13"
$ws.Cells.Item(5, 3).Value = 0.9304681353507841
$ws.Cells.Item(5, 4).Value = 0.9543142967306345
$ws.Cells.Item(6, 3).Value = 0.8867182122118554
$ws.Cells.Item(6, 4).Value = 0.9066219977820681
$ws.Cells.Item(6, 6).Value = [double]"1.354412011025889e-13"
$ws.Cells.Item(7, 3).Value = 0.8391446026191456
$ws.Cells.Item(7, 4).Value = 0.8668144507331865
$ws.Cells.Item(7, 6).Value = [double]"1.899434767816931e-13"
$ws.Cells.Item(8, 3).Value = 0.7928317718033387
$ws.Cells.Item(8, 4).Value = 0.8114747781125126
$ws.Cells.Item(8, 6).Value = [double]"2.275811002987144e-13"
$ws.Cells.Item(9, 3).Value = 0.7607776126703595
$ws.Cells.Item(9, 4).Value = 0.7741888599260756
$ws.Cells.Item(9, 6).Value = [double]"3.745620174033873e-13"
$ws.Cells.Item(10, 3).Value = 0.7326739465384652
$ws.Cells.Item(10, 4).Value = 0.7473663881730492
$ws.Cells.Item(10, 6).Value = [double]"3.139433694767128e-13"
$ws.Cells.Item(11, 3).Value = 0.7025972164436267
$ws.Cells.Item(11, 4).Value = 0.717981527635466
$ws.Cells.Item(11, 6).Value = [double]"2.943717855679378e-13"
$ws.Cells.Item(12, 3).Value = 0.6721163486365607
$ws.Cells.Item(12, 4).Value = 0.6872129288068618
$ws.Cells.Item(12, 6).Value = [double]"3.187850862341234e-13"
$ws.Cells.Item(13, 3).Value = 0.6437276347212123
$ws.Cells.Item(13, 4).Value = 0.6570197922644265
$ws.Cells.Item(13, 6).Value = [double]"3.475441159866836e-13"
$ws.Cells.Item(14, 3).Value = 0.6193606523952206
$ws.Cells.Item(14, 4).Value = 0.6304355006729732
$ws.Cells.Item(14, 6).Value = [double]"3.202706894299775e-13"
$ws.Cells.Item(15, 3).Value = 0.5988632757310339
$ws.Cells.Item(15, 4).Value = 0.6082858277787302
$ws.Cells.Item(15, 6).Value = [double]"2.835522585371362e-13"
$ws.Cells.Item(16, 3).Value = 0.5807429553881358
$ws.Cells.Item(16, 4).Value = 0.5894407477746866
$ws.Cells.Item(16, 6).Value = [double]"2.466320605509448e-13"
$ws.Cells.Item(17, 3).Value = 0.5631387433020947
$ws.Cells.Item(17, 4).Value = 0.5720451860200731
$ws.Cells.Item(17, 6).Value = [double]"2.180664924659718e-13"
$ws.Cells.Item(18, 3).Value = 0.5438958214577525
$ws.Cells.Item(18, 4).Value = 0.5542323237481271
$ws.Cells.Item(18, 6).Value = [double]"2.004109941404972e-13"
$ws.Cells.Item(19, 3).Value = 0.5206794560569086
$ws.Cells.Item(19, 4).Value = 0.5335593429547203
$ws.Cells.Item(19, 6).Value = [double]"1.730335697204537e-13"
$ws.Cells.Item(20, 3).Value = 0.4740770199242028
$ws.Cells.Item(20, 4).Value = 0.5077995931320249
$ws.Cells.Item(20, 6).Value = [double]"1.703961310839761e-13"
$ws.Cells.Item(21, 3).Value = 0.4269027628961122
$ws.Cells.Item(21, 4).Value = 0.4403544670642062
$ws.Cells.Item(21, 6).Value = [double]"2.825040073664013e-13"
$ws.Cells.Item(22, 3).Value = 0.4034390973046975
$ws.Cells.Item(22, 4).Value = 0.4134510778556927
$ws.Cells.Item(22, 6).Value = [double]"2.380278948604979e-13"
$ws.Cells.Item(23, 3).Value = 0.3854256130262844
$ws.Cells.Item(23, 4).Value = 0.3934271355826724
$ws.Cells.Item(23, 6).Value = [double]"2.020451493985933e-13"
$ws.Cells.Item(24, 3).Value = 0.37037748200145
$ws.Cells.Item(24, 4).Value = 0.3774241090637914
$ws.Cells.Item(24, 6).Value = [double]"1.663348517596428e-13"
$ws.Cells.Item(25, 3).Value = 0.3562569693507968
$ws.Cells.Item(25, 4).Value = 0.3633308736154864
$ws.Cells.Item(25, 6).Value = [double]"1.525841309939344e-13"
$ws.Cells.Item(26, 3).Value = 0.3414657465619925
$ws.Cells.Item(26, 4).Value = 0.3491830837318486
$ws.Cells.Item(26, 6).Value = [double]"1.454490851260519e-13"
$ws.Cells.Item(27, 3).Value = 0.315232606070213
$ws.Cells.Item(27, 4).Value = 0.3337484283854826
$ws.Cells.Item(27, 6).Value = [double]"5.374811491010176e-14"
